$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "MG745672"
$ws.Range("B5").Value = "FinchAvePV"
$ws.Range("C5").Value = "Pileated finch aveparvovirus"
$ws.Range("D5").Value = "Coryphospingus pileatus"
$ws.Range("E5").Value = "NK"
$ws.Range("F5").Value = "Aveparvovirus"
$ws.Range("G5").Value = "Parvovirinae"
$ws.Range("H5").Value = "NK"
$ws.Range("I5").Value = "NK"
$ws.Range("J5").Value = "NK"
$ws.Range("K5").Value = "NK"
$ws.Range("L5").Value = "NK"
$ws.Range("M5").Value = $false

$ws.Range("A1:M5").Select()
